$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" - append new row 79 ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Cells.Item(79,1).NumberFormat = $ws1.Cells.Item(78,1).NumberFormat
$ws1.Cells.Item(79,1).Value = 45669.99999999999
$ws1.Cells.Item(79,2).Value = 30

# --- Sheet 2: "Monthly Trend" - append new row 24 ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Cells.Item(24,1).NumberFormat = $ws2.Cells.Item(23,1).NumberFormat
$ws2.Cells.Item(24,1).Value = 45688.99999999999
$ws2.Cells.Item(24,2).Value = 30

# --- Sheet 3: "PO Forecast" - refreshed forecast model ---
$ws3 = $wb.Worksheets.Item("PO Forecast")

# Updated B-values for existing rows 2-78 (row 29 unchanged)
$ws3.Cells.Item(2,2).Value = 214
$ws3.Cells.Item(3,2).Value = 235
$ws3.Cells.Item(4,2).Value = 256
$ws3.Cells.Item(5,2).Value = 277
$ws3.Cells.Item(6,2).Value = 320
$ws3.Cells.Item(7,2).Value = 362
$ws3.Cells.Item(8,2).Value = 383
$ws3.Cells.Item(9,2).Value = 404
$ws3.Cells.Item(10,2).Value = 447
$ws3.Cells.Item(11,2).Value = 468
$ws3.Cells.Item(12,2).Value = 489
$ws3.Cells.Item(13,2).Value = 510
$ws3.Cells.Item(14,2).Value = 532
$ws3.Cells.Item(15,2).Value = 553
$ws3.Cells.Item(16,2).Value = 574
$ws3.Cells.Item(17,2).Value = 595
$ws3.Cells.Item(18,2).Value = 659
$ws3.Cells.Item(19,2).Value = 680
$ws3.Cells.Item(20,2).Value = 701
$ws3.Cells.Item(21,2).Value = 722
$ws3.Cells.Item(22,2).Value = 744
$ws3.Cells.Item(23,2).Value = 765
$ws3.Cells.Item(24,2).Value = 786
$ws3.Cells.Item(25,2).Value = 828
$ws3.Cells.Item(26,2).Value = 871
$ws3.Cells.Item(27,2).Value = 892
$ws3.Cells.Item(28,2).Value = 934
$ws3.Cells.Item(30,2).Value = 998
$ws3.Cells.Item(31,2).Value = 1019
$ws3.Cells.Item(32,2).Value = 1040
$ws3.Cells.Item(33,2).Value = 1062
$ws3.Cells.Item(34,2).Value = 1146
$ws3.Cells.Item(35,2).Value = 1168
$ws3.Cells.Item(36,2).Value = 1189
$ws3.Cells.Item(37,2).Value = 1210
$ws3.Cells.Item(38,2).Value = 1231
$ws3.Cells.Item(39,2).Value = 1295
$ws3.Cells.Item(40,2).Value = 1316
$ws3.Cells.Item(41,2).Value = 1337
$ws3.Cells.Item(42,2).Value = 1358
$ws3.Cells.Item(43,2).Value = 1380
$ws3.Cells.Item(44,2).Value = 1401
$ws3.Cells.Item(45,2).Value = 1422
$ws3.Cells.Item(46,2).Value = 1443
$ws3.Cells.Item(47,2).Value = 1464
$ws3.Cells.Item(48,2).Value = 1486
$ws3.Cells.Item(49,2).Value = 1507
$ws3.Cells.Item(50,2).Value = 1528
$ws3.Cells.Item(51,2).Value = 1549
$ws3.Cells.Item(52,2).Value = 1570
$ws3.Cells.Item(53,2).Value = 1592
$ws3.Cells.Item(54,2).Value = 1613
$ws3.Cells.Item(55,2).Value = 1634
$ws3.Cells.Item(56,2).Value = 1655
$ws3.Cells.Item(57,2).Value = 1676
$ws3.Cells.Item(58,2).Value = 1698
$ws3.Cells.Item(59,2).Value = 1719
$ws3.Cells.Item(60,2).Value = 1740
$ws3.Cells.Item(61,2).Value = 1761
$ws3.Cells.Item(62,2).Value = 1782
$ws3.Cells.Item(63,2).Value = 1804
$ws3.Cells.Item(64,2).Value = 1825
$ws3.Cells.Item(65,2).Value = 1867
$ws3.Cells.Item(66,2).Value = 1888
$ws3.Cells.Item(67,2).Value = 1910
$ws3.Cells.Item(68,2).Value = 1931
$ws3.Cells.Item(69,2).Value = 1952
$ws3.Cells.Item(70,2).Value = 1973
$ws3.Cells.Item(71,2).Value = 1994
$ws3.Cells.Item(72,2).Value = 2058
$ws3.Cells.Item(73,2).Value = 2079
$ws3.Cells.Item(74,2).Value = 2100
$ws3.Cells.Item(75,2).Value = 2122
$ws3.Cells.Item(76,2).Value = 2143
$ws3.Cells.Item(77,2).Value = 2164
$ws3.Cells.Item(78,2).Value = 2206

# New forecast tail rows 79-87 (dates + values shift forward)
$ws3.Cells.Item(79,1).NumberFormat = $ws3.Cells.Item(78,1).NumberFormat
$ws3.Cells.Item(79,1).Value = 45669.99999999999
$ws3.Cells.Item(79,2).Value = 2418
$ws3.Cells.Item(80,1).NumberFormat = $ws3.Cells.Item(79,1).NumberFormat
$ws3.Cells.Item(80,1).Value = 45676.99999999999
$ws3.Cells.Item(80,2).Value = 2440
$ws3.Cells.Item(81,1).NumberFormat = $ws3.Cells.Item(80,1).NumberFormat
$ws3.Cells.Item(81,1).Value = 45683.99999999999
$ws3.Cells.Item(81,2).Value = 2461
$ws3.Cells.Item(82,1).NumberFormat = $ws3.Cells.Item(81,1).NumberFormat
$ws3.Cells.Item(82,1).Value = 45690.99999999999
$ws3.Cells.Item(82,2).Value = 2482
$ws3.Cells.Item(83,1).NumberFormat = $ws3.Cells.Item(82,1).NumberFormat
$ws3.Cells.Item(83,1).Value = 45697.99999999999
$ws3.Cells.Item(83,2).Value = 2503
$ws3.Cells.Item(84,1).NumberFormat = $ws3.Cells.Item(83,1).NumberFormat
$ws3.Cells.Item(84,1).Value = 45704.99999999999
$ws3.Cells.Item(84,2).Value = 2524
$ws3.Cells.Item(85,1).NumberFormat = $ws3.Cells.Item(84,1).NumberFormat
$ws3.Cells.Item(85,1).Value = 45711.99999999999
$ws3.Cells.Item(85,2).Value = 2546
$ws3.Cells.Item(86,1).NumberFormat = $ws3.Cells.Item(85,1).NumberFormat
$ws3.Cells.Item(86,1).Value = 45718.99999999999
$ws3.Cells.Item(86,2).Value = 2567
$ws3.Cells.Item(87,1).NumberFormat = $ws3.Cells.Item(86,1).NumberFormat
$ws3.Cells.Item(87,1).Value = 45725.99999999999
$ws3.Cells.Item(87,2).Value = 2588
